# Atualizado por script em 02-11-2023 14:45
#
# 1) Rows 3-5 (F:V, the match/odds columns) got re-ordered: the match that
#    used to sit in row 4 now sits in row 3, the one in row 5 now sits in
#    row 4, and the one that used to be in row 3 now sits in row 5. Columns
#    A:E (Indice/pais/torneio/temporada/data_partida) stay put.
# 2) Four new match rows (57-60) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Cyclic rotation of rows 3, 4, 5 (columns F:V only) ----------------
$old3 = $ws.Range("F3:V3").Value()
$old4 = $ws.Range("F4:V4").Value()
$old5 = $ws.Range("F5:V5").Value()

$ws.Range("F3:V3").Value = $old4
$ws.Range("F4:V4").Value = $old5
$ws.Range("F5:V5").Value = $old3

# --- 2) Append new rows 57-60 ----------------------------------------------
$lastRow = 56

$newData = @(
    @(56, "iran", "persian-gulf-pro-league", "2023-2024", 45232.52083333334, "Aluminium Arak", 1, "Zob Ahan", 1, 2.44, "01/11/2023 00:42", 3.16, "02/11/2023 12:28", 2.54, "01/11/2023 00:42", 2.32, "02/11/2023 12:28", 3.21, "01/11/2023 00:42", 3.09, "02/11/2023 12:28", "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-zob-ahan/IV7qf6OS/"),
    @(57, "iran", "persian-gulf-pro-league", "2023-2024", 45232.54166666666, "Shams Azar Qazvin", 2, "Esteghlal F.C.", 2, 5.02, "01/11/2023 00:42", 5.08, "02/11/2023 12:59", 3.4, "01/11/2023 00:42", 2.99, "02/11/2023 12:59", 1.63, "01/11/2023 00:42", 1.85, "02/11/2023 12:59", "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/shams-azar-qazvin-esteghlal-teh/StFWd8vA/"),
    @(58, "iran", "persian-gulf-pro-league", "2023-2024", 45232.54166666666, "Tractor", 3, "Nassaji Mazandaran", 0, 1.58, "01/11/2023 01:12", 1.71, "02/11/2023 12:58", 3.35, "01/11/2023 01:12", 3.26, "02/11/2023 12:58", 5.39, "01/11/2023 01:12", 5.53, "02/11/2023 12:58", "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/tractor-mazandaran/lEbFlpnj/"),
    @(59, "iran", "persian-gulf-pro-league", "2023-2024", 45232.58333333334, "Sepahan", 2, "Malavan", 3, 1.34, "01/11/2023 02:12", 1.44, "02/11/2023 13:57", 4.27, "01/11/2023 02:12", 4.03, "02/11/2023 13:58", 7.52, "01/11/2023 02:12", 7.75, "02/11/2023 13:58", "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sepahan-malavan/S4cJm41d/")
)

$newRowsCount = $newData.Count

# Clone the formatting (bold index style on column A, date style on column E)
# from the last existing row onto each new row.
for ($i = 1; $i -le $newRowsCount; $i++) {
    $target = $lastRow + $i
    $ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
    $ws.Range("A" + $target + ":V" + $target).PasteSpecial(-4122)
}

# PowerShell's @(@(...)) literal builds a jagged System.Object[], not a true
# rectangular System.Object[,] - the COM bridge here silently ignores
# multi-cell writes that aren't a genuine 2D array, so build one explicitly.
$nCols = $newData[0].Count
$arr = New-Object 'object[,]' $newRowsCount,$nCols
for ($r = 0; $r -lt $newRowsCount; $r++) {
    for ($c = 0; $c -lt $nCols; $c++) {
        $arr[$r,$c] = $newData[$r][$c]
    }
}

$ws.Range("A" + ($lastRow + 1) + ":V" + ($lastRow + $newRowsCount)).Value = $arr
